# --- Scheduled-runner data refresh for Raiden_Profits --------------------
# Recomputed currentAveragePrice / LevePrice / LeveProfit figures (columns
# H-N) for a batch of leves across the ALC, ARM, BSM, CRP, GSM, LTW and WVR
# crafting sheets, pulling in updated market-board averages. Two stale cells
# (GSM!N5 and WVR!M137) are cleared outright because their column no longer
# carries a value for that row after the refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1521.8
$ws.Cells.Item(17, 10).Value = 1521.8
$ws.Cells.Item(17, 12).Value = 4565.4
$ws.Cells.Item(17, 14).Value = -4901.4
$ws.Cells.Item(76, 8).Value = 7566.5713
$ws.Cells.Item(76, 9).Value = 5747.5
$ws.Cells.Item(76, 11).Value = 5747.5
$ws.Cells.Item(76, 13).Value = -5432.5
$ws.Cells.Item(79, 8).Value = 7566.5713
$ws.Cells.Item(79, 9).Value = 5747.5
$ws.Cells.Item(79, 11).Value = 5747.5
$ws.Cells.Item(79, 13).Value = -4655.5
$ws.Cells.Item(86, 8).Value = 1785.5
$ws.Cells.Item(86, 9).Value = 1477.4
$ws.Cells.Item(86, 10).Value = 2093.6
$ws.Cells.Item(86, 11).Value = 1477.4
$ws.Cells.Item(86, 12).Value = 2093.6
$ws.Cells.Item(86, 13).Value = -354.4000000000001
$ws.Cells.Item(86, 14).Value = -4339.6
$ws.Cells.Item(89, 8).Value = 1785.5
$ws.Cells.Item(89, 9).Value = 1477.4
$ws.Cells.Item(89, 10).Value = 2093.6
$ws.Cells.Item(89, 11).Value = 7387.0
$ws.Cells.Item(89, 12).Value = 10468.0
$ws.Cells.Item(89, 13).Value = -1771.0
$ws.Cells.Item(89, 14).Value = -21700.0
$ws.Cells.Item(112, 8).Value = 3178.8076
$ws.Cells.Item(112, 10).Value = 3253.96
$ws.Cells.Item(112, 12).Value = 9761.880000000001
$ws.Cells.Item(112, 14).Value = -11977.88
$ws.Cells.Item(113, 8).Value = 7027.5386
$ws.Cells.Item(113, 9).Value = 7078.636
$ws.Cells.Item(113, 11).Value = 7078.636
$ws.Cells.Item(113, 13).Value = -3824.636
$ws.Cells.Item(135, 8).Value = 2475.8572
$ws.Cells.Item(135, 10).Value = 1800.0
$ws.Cells.Item(135, 12).Value = 16200.0
$ws.Cells.Item(135, 14).Value = -21270.0
$ws.Cells.Item(137, 8).Value = 3374.1
$ws.Cells.Item(137, 9).Value = 2344.842
$ws.Cells.Item(137, 11).Value = 7034.526
$ws.Cells.Item(137, 13).Value = -4484.526
$ws.Cells.Item(138, 8).Value = 2798.957
$ws.Cells.Item(138, 10).Value = 2946.9138
$ws.Cells.Item(138, 12).Value = 8840.741399999999
$ws.Cells.Item(138, 14).Value = -19120.7414
$ws.Cells.Item(141, 8).Value = 5222.1763
$ws.Cells.Item(141, 9).Value = 4481.4165
$ws.Cells.Item(141, 10).Value = 7000.0
$ws.Cells.Item(141, 11).Value = 13444.2495
$ws.Cells.Item(141, 12).Value = 21000.0
$ws.Cells.Item(141, 13).Value = -8264.249500000002
$ws.Cells.Item(141, 14).Value = -31360.0

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3511.2273
$ws.Cells.Item(61, 10).Value = 7989.0
$ws.Cells.Item(61, 12).Value = 7989.0
$ws.Cells.Item(61, 14).Value = -8413.0
$ws.Cells.Item(136, 8).Value = 3511.2273
$ws.Cells.Item(136, 10).Value = 7989.0
$ws.Cells.Item(136, 12).Value = 23967.0
$ws.Cells.Item(136, 14).Value = -29067.0

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 911.9
$ws.Cells.Item(22, 9).Value = 760.0
$ws.Cells.Item(22, 10).Value = 1013.1667
$ws.Cells.Item(22, 11).Value = 760.0
$ws.Cells.Item(22, 12).Value = 1013.1667
$ws.Cells.Item(22, 13).Value = -587.0
$ws.Cells.Item(22, 14).Value = -1359.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 6316.8184
$ws.Cells.Item(4, 9).Value = 7176.4287
$ws.Cells.Item(4, 10).Value = 4812.5
$ws.Cells.Item(4, 11).Value = 7176.4287
$ws.Cells.Item(4, 12).Value = 4812.5
$ws.Cells.Item(4, 13).Value = -7064.4287
$ws.Cells.Item(4, 14).Value = -5036.5
$ws.Cells.Item(16, 8).Value = 1249.1538
$ws.Cells.Item(16, 9).Value = 1185.3636
$ws.Cells.Item(16, 10).Value = 1600.0
$ws.Cells.Item(16, 11).Value = 1185.3636
$ws.Cells.Item(16, 12).Value = 1600.0
$ws.Cells.Item(16, 13).Value = -898.3635999999999
$ws.Cells.Item(16, 14).Value = -2174.0
$ws.Cells.Item(31, 8).Value = 7390.614
$ws.Cells.Item(31, 9).Value = 4128.2144
$ws.Cells.Item(31, 10).Value = 13099.8125
$ws.Cells.Item(31, 11).Value = 4128.2144
$ws.Cells.Item(31, 12).Value = 13099.8125
$ws.Cells.Item(31, 13).Value = -3833.2144
$ws.Cells.Item(31, 14).Value = -13689.8125
$ws.Cells.Item(34, 8).Value = 7390.614
$ws.Cells.Item(34, 9).Value = 4128.2144
$ws.Cells.Item(34, 10).Value = 13099.8125
$ws.Cells.Item(34, 11).Value = 4128.2144
$ws.Cells.Item(34, 12).Value = 13099.8125
$ws.Cells.Item(34, 13).Value = -3926.2144
$ws.Cells.Item(34, 14).Value = -13503.8125
$ws.Cells.Item(113, 8).Value = 1249.1538
$ws.Cells.Item(113, 9).Value = 1185.3636
$ws.Cells.Item(113, 10).Value = 1600.0
$ws.Cells.Item(113, 11).Value = 1185.3636
$ws.Cells.Item(113, 12).Value = 1600.0
$ws.Cells.Item(113, 13).Value = 984.6364000000001
$ws.Cells.Item(113, 14).Value = -5940.0
$ws.Cells.Item(122, 8).Value = 2267.238
$ws.Cells.Item(122, 9).Value = 2280.6
$ws.Cells.Item(122, 10).Value = 2000.0
$ws.Cells.Item(122, 11).Value = 6841.799999999999
$ws.Cells.Item(122, 12).Value = 6000.0
$ws.Cells.Item(122, 13).Value = -4391.799999999999
$ws.Cells.Item(122, 14).Value = -10900.0

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 10000.0
$ws.Cells.Item(5, 9).Value = 10000.0
$ws.Cells.Item(5, 10).Value = 0.0
$ws.Cells.Item(5, 11).Value = 10000.0
$ws.Cells.Item(5, 12).Value = 0.0
$ws.Cells.Item(5, 13).Value = -9888.0
$ws.Cells.Item(5, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 2955.4695
$ws.Cells.Item(132, 9).Value = 3147.8333
$ws.Cells.Item(132, 11).Value = 9443.499899999999
$ws.Cells.Item(132, 13).Value = -6913.499899999999
$ws.Cells.Item(134, 8).Value = 105999.664
$ws.Cells.Item(134, 10).Value = 105999.664
$ws.Cells.Item(134, 12).Value = 317998.992
$ws.Cells.Item(134, 14).Value = -323068.992

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 670.0
$ws.Cells.Item(27, 8).Value = 670.0
$ws.Cells.Item(61, 8).Value = 2496.2
$ws.Cells.Item(61, 9).Value = 2408.0
$ws.Cells.Item(61, 10).Value = 2849.0
$ws.Cells.Item(61, 11).Value = 2408.0
$ws.Cells.Item(61, 12).Value = 2849.0
$ws.Cells.Item(61, 13).Value = -2206.0
$ws.Cells.Item(61, 14).Value = -3253.0
$ws.Cells.Item(82, 8).Value = 1949.4
$ws.Cells.Item(82, 9).Value = 1949.4
$ws.Cells.Item(82, 11).Value = 1949.4
$ws.Cells.Item(82, 13).Value = -1588.4
$ws.Cells.Item(85, 8).Value = 1949.4
$ws.Cells.Item(85, 9).Value = 1949.4
$ws.Cells.Item(85, 11).Value = 1949.4
$ws.Cells.Item(85, 13).Value = -701.4000000000001
$ws.Cells.Item(113, 8).Value = 2496.2
$ws.Cells.Item(113, 9).Value = 2408.0
$ws.Cells.Item(113, 10).Value = 2849.0
$ws.Cells.Item(113, 11).Value = 2408.0
$ws.Cells.Item(113, 12).Value = 2849.0
$ws.Cells.Item(113, 13).Value = -238.0
$ws.Cells.Item(113, 14).Value = -7189.0
$ws.Cells.Item(132, 8).Value = 3218.524
$ws.Cells.Item(132, 9).Value = 3339.0
$ws.Cells.Item(132, 11).Value = 10017.0
$ws.Cells.Item(132, 13).Value = -7487.0
$ws.Cells.Item(136, 8).Value = 3769.8857
$ws.Cells.Item(136, 9).Value = 3606.1667
$ws.Cells.Item(136, 11).Value = 10818.5001
$ws.Cells.Item(136, 13).Value = -8268.500100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 4296.5
$ws.Cells.Item(81, 9).Value = 2399.0
$ws.Cells.Item(81, 11).Value = 4798.0
$ws.Cells.Item(81, 13).Value = -3737.0
$ws.Cells.Item(84, 8).Value = 4296.5
$ws.Cells.Item(84, 9).Value = 2399.0
$ws.Cells.Item(84, 11).Value = 23990.0
$ws.Cells.Item(84, 13).Value = -18686.0
$ws.Cells.Item(132, 8).Value = 2296.2693
$ws.Cells.Item(132, 9).Value = 2030.7391
$ws.Cells.Item(132, 11).Value = 6092.2173
$ws.Cells.Item(132, 13).Value = -3562.2173
$ws.Cells.Item(136, 8).Value = 3064.8948
$ws.Cells.Item(136, 9).Value = 3024.375
$ws.Cells.Item(136, 11).Value = 9073.125
$ws.Cells.Item(136, 13).Value = -6523.125
$ws.Cells.Item(137, 8).Value = 0.0
$ws.Cells.Item(137, 9).Value = 0.0
$ws.Cells.Item(137, 11).Value = 0.0
$ws.Cells.Item(137, 13).ClearContents()
